$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new "Elbo Room Bar" webcam entries (Fort Lauderdale, FL, USA)
# below the existing data, continuing the table in rows 213-216.

$data = @(
    @("LIVE, PUB, CLUB", "26.119447755690356, -80.10472547327616", "Elbo Room Bar WebCam 1", "FL", "USA", "KY4Yd5QR570"),
    @("LIVE, PUB, CLUB", "26.119469408514355, -80.10463347190057", "Elbo Room Bar WebCam 2", "FL", "USA", "yqx3qFTGVEY"),
    @("LIVE, PUB, CLUB", "26.119456764955224, -80.1046965038119",  "Elbo Room Bar WebCam 3", "FL", "USA", "ZGnLVjkBEls"),
    @("LIVE, PUB, CLUB", "26.119225568236672, -80.10454294734677", "Elbo Room Bar WebCam 4", "FL", "USA", "uoJRuZg3NME")
)

$startRow = 213
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]

    # Carry the bordered look of the row above down into the new rows
    # (Category and Country columns are bordered throughout the sheet).
    $ws.Range("A212").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Range("E212").Copy()
    $ws.Cells.Item($r, 5).PasteSpecial(-4122)
}

# Leave the selection where the user last left off, on the final new row.
$ws.Range("D216").Select()
